$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that were removed in the diff ---
$ws.Range("B4").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("A42").ClearContents()
$ws.Range("B47").ClearContents()
$ws.Range("B52").ClearContents()
$ws.Range("B55").ClearContents()
$ws.Range("B101").ClearContents()

# --- Update cells whose values changed ---
$ws.Range("B71").Value = 855800
$ws.Range("B80").Value = 52822
$ws.Range("B83").Value = 855000
$ws.Range("B85").Value = 84000
$ws.Range("B94").Value = 73000

# --- Update the view / selection state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$ws.Range("A42").Select() | Out-Null
